# The "FORBUDES 400/12MCG 60 INHALATION CAPS.+INHALER" line item (row 29)
# was removed from the day-sale report. Deleting the whole row shifts every
# row below it up by one, which also naturally fixes the merged-cell ranges
# further down the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(29).Delete()

# The running total at the bottom of the report (now on row 68, column P)
# needs to drop by the removed item's price (334.00).
$ws.Cells.Item(68, 16).Value = 3439.5949999999998

# The report was re-generated/re-uploaded a minute later, so the footer
# timestamp (now on row 69, column A) moves from 2:53 PM to 2:54 PM.
$ws.Cells.Item(69, 1).Value = "Monday, 21 July, 2025 2:54 PM"
